# Update cryptocurrency price/volume data per commit:
# "Updated cryptos list on Sun Jan 21 08:39:20 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price (D) and Volume (E) columns are stored as text in this sheet,
# so force text number format before assigning values to avoid Excel
# auto-converting numeric-looking strings (e.g. "15.70") into numbers
# (which would drop formatting like trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.676.25"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.475.82"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.49"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.79"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.16%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.11"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +8.67%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.857.56"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.91"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.70"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.490.87"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.627.74"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0950"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.40"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.31"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.67"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.74"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.93"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.77"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.85"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.25"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.34"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.52"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.38"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.93"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.85"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.05%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.987.48"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0285"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.93"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.30"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.715.27"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.48"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.23"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.26"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.18%  "
